$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate column A with the home-page menu labels (mirrors the authored
# "home page style" edit: shared-string backed rows A1:A8).
$ws.Range("A1").Value = "运行日志"
$ws.Range("A2").Value = "电量统计"
$ws.Range("A3").Value = "库房"
$ws.Range("A4").Value = "案例库"
$ws.Range("A5").Value = "两票"
$ws.Range("A6").Value = "设备缺陷"
$ws.Range("A7").Value = "档案"
$ws.Range("A8").Value = "考勤"

# Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it.
[void]$ws.Range("F11").Select()
